$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50 holds the "N" (sample size) figures which were halved as part of
# adding an extra loop layer for UK totals. Update B50:G50 from 94660 to 47330.
$ws.Range("B50:G50").Value = 47330
